# Issue #219 change to build Any CPU and refactor Excel and CSV reading code to support
# The header for the account number column was renamed from the display label
# "Account Number" to the schema/logical name "accountnumber".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "accountnumber"
